# MigrationRenamer.xlsx — rebase the migration file list onto a new
# timestamp (2021_03_25_) and swap the old "media" migration for a
# "settings" migration that already carries its own generated timestamp
# (2021_03_25_150256_create_settings_table.php). The "remote_items"
# migration is also renamed to "remote". All the other columns
# (Name/Date/Seq/Replace/CMD) are calculated Table1 columns that
# recompute automatically once the Files column is edited.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value  = "2021_03_25_000001_update_users_table.php"
$ws.Range("A3").Value  = "2021_03_25_000002_create_jobs_table.php"
$ws.Range("A4").Value  = "2021_03_25_150256_create_settings_table.php"
$ws.Range("A5").Value  = "2021_03_25_000004_create_masters_table.php"
$ws.Range("A6").Value  = "2021_03_25_000005_create_user_logins_table.php"
$ws.Range("A7").Value  = "2021_03_25_000006_create_items_table.php"
$ws.Range("A8").Value  = "2021_03_25_000007_create_item_groups_table.php"
$ws.Range("A9").Value  = "2021_03_25_000008_create_menus_table.php"
$ws.Range("A10").Value = "2021_03_25_000009_create_price_lists_table.php"
$ws.Range("A11").Value = "2021_03_25_000010_create_prices_table.php"
$ws.Range("A12").Value = "2021_03_25_000011_create_taxes_table.php"
$ws.Range("A13").Value = "2021_03_25_000012_create_kitchens_table.php"
$ws.Range("A14").Value = "2021_03_25_000013_create_kitchen_items_table.php"
$ws.Range("A15").Value = "2021_03_25_000014_create_kitchen_statuses_table.php"
$ws.Range("A16").Value = "2021_03_25_000015_create_customers_table.php"
$ws.Range("A17").Value = "2021_03_25_000016_create_seatings_table.php"
$ws.Range("A18").Value = "2021_03_25_000017_create_tokens_table.php"
$ws.Range("A19").Value = "2021_03_25_000018_create_token_items_table.php"
$ws.Range("A20").Value = "2021_03_25_000019_create_bills_table.php"
$ws.Range("A21").Value = "2021_03_25_000020_create_payments_table.php"
$ws.Range("A22").Value = "2021_03_25_000021_create_remote_table.php"

# The "Date" calculated column literally bakes in today's migration
# stamp, so the column formula + header formula text move forward too.
$ws.ListObjects.Item("Table1").ListColumns.Item("Date").DataBodyRange.Formula = '="2021_03_25_"'

# Scroll the view over one column, matching the saved selection/viewport.
$ws.Range("B1").Select()

$excel.Calculate()
